# Update New Orleans xlsx: reorder sheets (review_info first, hotel_info second)
# and add a new "State" column to hotel_info (inserted between Hotel_Name and City)
# with the value "Louisiana" for the existing hotel row.

$wb = $excel.ActiveWorkbook

$hotel = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# Insert a new column C in hotel_info, between Hotel_Name (B) and City (old C)
# (shifts old City, Zip, ... one column right)
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# Reorder worksheet tabs so review_info comes before hotel_info
$review.Move($hotel)
